# Natmi following Dr Hou advice
# Rebuild the LR-pairs (Col4a3-Itgav) sheet so that "M2" is included as a
# sending cluster as well as a target cluster (previously only ECs/FAPs/sCs
# were senders). This turns the 3x4 sender/target grid into a 4x4 grid
# (rows 2-17) and refreshes every downstream statistic accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Col4a3"
$ws.Cells.Item(2, 3).Value = "Itgav"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.6813643333333333
$ws.Cells.Item(2, 8).Value = 2.044093
$ws.Cells.Item(2, 9).Value = 0.3450055461675409
$ws.Cells.Item(2, 10).Value = 0.3450055461675409
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 13.441269
$ws.Cells.Item(2, 14).Value = 40.323807
$ws.Cells.Item(2, 15).Value = 0.0897308213348123
$ws.Cells.Item(2, 16).Value = 0.08973082133481232
$ws.Cells.Item(2, 17).Value = 9.158401291339
$ws.Cells.Item(2, 18).Value = 82.42561162205101
$ws.Cells.Item(2, 19).Value = 0.03095763102267895
$ws.Cells.Item(2, 20).Value = 0.03095763102267896

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Col4a3"
$ws.Cells.Item(3, 3).Value = "Itgav"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.6813643333333333
$ws.Cells.Item(3, 8).Value = 2.044093
$ws.Cells.Item(3, 9).Value = 0.3450055461675409
$ws.Cells.Item(3, 10).Value = 0.3450055461675409
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 54.711535
$ws.Cells.Item(3, 14).Value = 164.134605
$ws.Cells.Item(3, 15).Value = 0.3652416280068742
$ws.Cells.Item(3, 16).Value = 0.3652416280068742
$ws.Cells.Item(3, 17).Value = 37.27848857091833
$ws.Cells.Item(3, 18).Value = 335.506397138265
$ws.Cells.Item(3, 19).Value = 0.1260103873536334
$ws.Cells.Item(3, 20).Value = 0.1260103873536335

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Col4a3"
$ws.Cells.Item(4, 3).Value = "Itgav"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.6813643333333333
$ws.Cells.Item(4, 8).Value = 2.044093
$ws.Cells.Item(4, 9).Value = 0.3450055461675409
$ws.Cells.Item(4, 10).Value = 0.3450055461675409
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 63.67711
$ws.Cells.Item(4, 14).Value = 191.03133
$ws.Cells.Item(4, 15).Value = 0.4250937452800914
$ws.Cells.Item(4, 16).Value = 0.4250937452800915
$ws.Cells.Item(4, 17).Value = 43.38731160374333
$ws.Cells.Item(4, 18).Value = 390.48580443369
$ws.Cells.Item(4, 19).Value = 0.1466596997627634
$ws.Cells.Item(4, 20).Value = 0.1466596997627635

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Col4a3"
$ws.Cells.Item(5, 3).Value = "Itgav"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.6813643333333333
$ws.Cells.Item(5, 8).Value = 2.044093
$ws.Cells.Item(5, 9).Value = 0.3450055461675409
$ws.Cells.Item(5, 10).Value = 0.3450055461675409
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 17.96553866666667
$ws.Cells.Item(5, 14).Value = 53.896616
$ws.Cells.Item(5, 15).Value = 0.119933805378222
$ws.Cells.Item(5, 16).Value = 0.119933805378222
$ws.Cells.Item(5, 17).Value = 12.24107727658756
$ws.Cells.Item(5, 18).Value = 110.169695489288
$ws.Cells.Item(5, 19).Value = 0.04137782802846503
$ws.Cells.Item(5, 20).Value = 0.04137782802846504

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Col4a3"
$ws.Cells.Item(6, 3).Value = "Itgav"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.845217
$ws.Cells.Item(6, 8).Value = 2.535651
$ws.Cells.Item(6, 9).Value = 0.4279715542029013
$ws.Cells.Item(6, 10).Value = 0.4279715542029013
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 13.441269
$ws.Cells.Item(6, 14).Value = 40.323807
$ws.Cells.Item(6, 15).Value = 0.0897308213348123
$ws.Cells.Item(6, 16).Value = 0.08973082133481232
$ws.Cells.Item(6, 17).Value = 11.360789060373
$ws.Cells.Item(6, 18).Value = 102.247101543357
$ws.Cells.Item(6, 19).Value = 0.03840223906656248
$ws.Cells.Item(6, 20).Value = 0.03840223906656249

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Col4a3"
$ws.Cells.Item(7, 3).Value = "Itgav"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.845217
$ws.Cells.Item(7, 8).Value = 2.535651
$ws.Cells.Item(7, 9).Value = 0.4279715542029013
$ws.Cells.Item(7, 10).Value = 0.4279715542029013
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 54.711535
$ws.Cells.Item(7, 14).Value = 164.134605
$ws.Cells.Item(7, 15).Value = 0.3652416280068742
$ws.Cells.Item(7, 16).Value = 0.3652416280068742
$ws.Cells.Item(7, 17).Value = 46.243119478095
$ws.Cells.Item(7, 18).Value = 416.188075302855
$ws.Cells.Item(7, 19).Value = 0.1563130271976998
$ws.Cells.Item(7, 20).Value = 0.1563130271976999

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Col4a3"
$ws.Cells.Item(8, 3).Value = "Itgav"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.845217
$ws.Cells.Item(8, 8).Value = 2.535651
$ws.Cells.Item(8, 9).Value = 0.4279715542029013
$ws.Cells.Item(8, 10).Value = 0.4279715542029013
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 63.67711
$ws.Cells.Item(8, 14).Value = 191.03133
$ws.Cells.Item(8, 15).Value = 0.4250937452800914
$ws.Cells.Item(8, 16).Value = 0.4250937452800915
$ws.Cells.Item(8, 17).Value = 53.82097588287
$ws.Cells.Item(8, 18).Value = 484.38878294583
$ws.Cells.Item(8, 19).Value = 0.1819280308494529
$ws.Cells.Item(8, 20).Value = 0.181928030849453

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Col4a3"
$ws.Cells.Item(9, 3).Value = "Itgav"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.845217
$ws.Cells.Item(9, 8).Value = 2.535651
$ws.Cells.Item(9, 9).Value = 0.4279715542029013
$ws.Cells.Item(9, 10).Value = 0.4279715542029013
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 17.96553866666667
$ws.Cells.Item(9, 14).Value = 53.896616
$ws.Cells.Item(9, 15).Value = 0.119933805378222
$ws.Cells.Item(9, 16).Value = 0.119933805378222
$ws.Cells.Item(9, 17).Value = 15.184778695224
$ws.Cells.Item(9, 18).Value = 136.663008257016
$ws.Cells.Item(9, 19).Value = 0.05132825708918595
$ws.Cells.Item(9, 20).Value = 0.05132825708918596

$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Col4a3"
$ws.Cells.Item(10, 3).Value = "Itgav"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.009069333333333334
$ws.Cells.Item(10, 8).Value = 0.027208
$ws.Cells.Item(10, 9).Value = 0.004592213221280271
$ws.Cells.Item(10, 10).Value = 0.004592213221280271
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 13.441269
$ws.Cells.Item(10, 14).Value = 40.323807
$ws.Cells.Item(10, 15).Value = 0.0897308213348123
$ws.Cells.Item(10, 16).Value = 0.08973082133481232
$ws.Cells.Item(10, 17).Value = 0.121903348984
$ws.Cells.Item(10, 18).Value = 1.097130140856
$ws.Cells.Item(10, 19).Value = 0.0004120630640900628
$ws.Cells.Item(10, 20).Value = 0.0004120630640900629

$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Col4a3"
$ws.Cells.Item(11, 3).Value = "Itgav"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.009069333333333334
$ws.Cells.Item(11, 8).Value = 0.027208
$ws.Cells.Item(11, 9).Value = 0.004592213221280271
$ws.Cells.Item(11, 10).Value = 0.004592213221280271
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 54.711535
$ws.Cells.Item(11, 14).Value = 164.134605
$ws.Cells.Item(11, 15).Value = 0.3652416280068742
$ws.Cells.Item(11, 16).Value = 0.3652416280068742
$ws.Cells.Item(11, 17).Value = 0.4961971480933333
$ws.Cells.Item(11, 18).Value = 4.46577433284
$ws.Cells.Item(11, 19).Value = 0.001677267433095098
$ws.Cells.Item(11, 20).Value = 0.001677267433095098

$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Col4a3"
$ws.Cells.Item(12, 3).Value = "Itgav"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.009069333333333334
$ws.Cells.Item(12, 8).Value = 0.027208
$ws.Cells.Item(12, 9).Value = 0.004592213221280271
$ws.Cells.Item(12, 10).Value = 0.004592213221280271
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 63.67711
$ws.Cells.Item(12, 14).Value = 191.03133
$ws.Cells.Item(12, 15).Value = 0.4250937452800914
$ws.Cells.Item(12, 16).Value = 0.4250937452800915
$ws.Cells.Item(12, 17).Value = 0.5775089362933333
$ws.Cells.Item(12, 18).Value = 5.19758042664
$ws.Cells.Item(12, 19).Value = 0.001952121117358783
$ws.Cells.Item(12, 20).Value = 0.001952121117358784

$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Col4a3"
$ws.Cells.Item(13, 3).Value = "Itgav"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.009069333333333334
$ws.Cells.Item(13, 8).Value = 0.027208
$ws.Cells.Item(13, 9).Value = 0.004592213221280271
$ws.Cells.Item(13, 10).Value = 0.004592213221280271
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 17.96553866666667
$ws.Cells.Item(13, 14).Value = 53.896616
$ws.Cells.Item(13, 15).Value = 0.119933805378222
$ws.Cells.Item(13, 16).Value = 0.119933805378222
$ws.Cells.Item(13, 17).Value = 0.1629354586808889
$ws.Cells.Item(13, 18).Value = 1.466419128128
$ws.Cells.Item(13, 19).Value = 0.0005507616067363259
$ws.Cells.Item(13, 20).Value = 0.000550761606736326

$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Col4a3"
$ws.Cells.Item(14, 3).Value = "Itgav"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.4392866666666667
$ws.Cells.Item(14, 8).Value = 1.31786
$ws.Cells.Item(14, 9).Value = 0.2224306864082776
$ws.Cells.Item(14, 10).Value = 0.2224306864082776
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 13.441269
$ws.Cells.Item(14, 14).Value = 40.323807
$ws.Cells.Item(14, 15).Value = 0.0897308213348123
$ws.Cells.Item(14, 16).Value = 0.08973082133481232
$ws.Cells.Item(14, 17).Value = 5.90457025478
$ws.Cells.Item(14, 18).Value = 53.14113229302001
$ws.Cells.Item(14, 19).Value = 0.01995888818148082
$ws.Cells.Item(14, 20).Value = 0.01995888818148083

$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Col4a3"
$ws.Cells.Item(15, 3).Value = "Itgav"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.4392866666666667
$ws.Cells.Item(15, 8).Value = 1.31786
$ws.Cells.Item(15, 9).Value = 0.2224306864082776
$ws.Cells.Item(15, 10).Value = 0.2224306864082776
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 54.711535
$ws.Cells.Item(15, 14).Value = 164.134605
$ws.Cells.Item(15, 15).Value = 0.3652416280068742
$ws.Cells.Item(15, 16).Value = 0.3652416280068742
$ws.Cells.Item(15, 17).Value = 24.03404783836666
$ws.Cells.Item(15, 18).Value = 216.3064305453
$ws.Cells.Item(15, 19).Value = 0.0812409460224458
$ws.Cells.Item(15, 20).Value = 0.08124094602244583

$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Col4a3"
$ws.Cells.Item(16, 3).Value = "Itgav"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.4392866666666667
$ws.Cells.Item(16, 8).Value = 1.31786
$ws.Cells.Item(16, 9).Value = 0.2224306864082776
$ws.Cells.Item(16, 10).Value = 0.2224306864082776
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 63.67711
$ws.Cells.Item(16, 14).Value = 191.03133
$ws.Cells.Item(16, 15).Value = 0.4250937452800914
$ws.Cells.Item(16, 16).Value = 0.4250937452800915
$ws.Cells.Item(16, 17).Value = 27.97250539486667
$ws.Cells.Item(16, 18).Value = 251.7525485538
$ws.Cells.Item(16, 19).Value = 0.09455389355051624
$ws.Cells.Item(16, 20).Value = 0.09455389355051627

$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Col4a3"
$ws.Cells.Item(17, 3).Value = "Itgav"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.4392866666666667
$ws.Cells.Item(17, 8).Value = 1.31786
$ws.Cells.Item(17, 9).Value = 0.2224306864082776
$ws.Cells.Item(17, 10).Value = 0.2224306864082776
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 17.96553866666667
$ws.Cells.Item(17, 14).Value = 53.896616
$ws.Cells.Item(17, 15).Value = 0.119933805378222
$ws.Cells.Item(17, 16).Value = 0.119933805378222
$ws.Cells.Item(17, 17).Value = 7.892021595751111
$ws.Cells.Item(17, 18).Value = 71.02819436176
$ws.Cells.Item(17, 19).Value = 0.0266769586538347
$ws.Cells.Item(17, 20).Value = 0.0266769586538347
